# Update the "Install package" syntax command (B11) to reflect the new
# repository location: the TomLeversRPackage folder now lives inside the
# "R" repository folder, so the install.packages() source path needs an
# extra "\R\" path segment inserted before "TomLeversRPackage".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSyntax = "install.packages(""tidyverse"", repos = ""http://cran.us.r-project.org"")`ninstall.packages(""C:\\Users\\Tom\\Documents\\Tom_Levers_Git_Repository\\R\\TomLeversRPackage"", repos = NULL, type=""source"")"

$ws.Range("B11").Value = $newSyntax

# Column B widened slightly to continue to best-fit the (now longer) text.
$ws.Columns.Item(2).ColumnWidth = 110.83
